# Generate Report for Handoff
# Adds a new handed-off file ("6c17df13-cd98-445e-a0f6-ca6b5f41dfe2.md") as a
# new row (row 9) on all three worksheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileName   = "6c17df13-cd98-445e-a0f6-ca6b5f41dfe2.md"
$pathName   = "e2e\6c17df13-cd98-445e-a0f6-ca6b5f41dfe2.md"
$commitHash = "f1635c37182a5c75cee6efd691fff00c118c1527"
$ghBase     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitHash/e2e/$fileName"

$xliffZhCn  = "6c17df13-cd98-445e-a0f6-ca6b5f41dfe2.f1635c37182a5c75cee6efd691fff00c118c1527.zh-cn.xlf"
$xliffDeDe  = "6c17df13-cd98-445e-a0f6-ca6b5f41dfe2.f1635c37182a5c75cee6efd691fff00c118c1527.de-de.xlf"

$handoffDateOverview = "2016-12-15 03:51:26"
$handoffDateZhCn     = "2016-12-15 03:51:13"
$handoffDateDeDe     = "2016-12-15 03:51:26"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | Path And Name | Extension | Publish URL |
#                 zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(9, 1).Value = $fileName
$wsOverview.Cells.Item(9, 2).Value = $pathName
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(9, 2), $ghBase, "", "", $pathName) | Out-Null
$wsOverview.Cells.Item(9, 3).Value = ".md"
$wsOverview.Cells.Item(9, 4).Value = ""
$wsOverview.Cells.Item(9, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item(9, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item(9, 7).Value = $handoffDateOverview

# ---------------------------------------------------------------------------
# zh-cn / de-de sheets share the same column layout:
# Source File Name | File Extension | Status | Source Path | Priority |
# Content Duplicate | Latest Handoff File | Latest Handoff Datetime |
# Lastest Handoff Name | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Latest Handback Name | Reference Tokens |
# To be localized | Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------------------

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Cells.Item(9, 1).Value = $fileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item(9, 1), $ghBase, "", "", $fileName) | Out-Null
$wsZhCn.Cells.Item(9, 2).Value = ".md"
$wsZhCn.Cells.Item(9, 3).Value = "Ready for handoff"
$wsZhCn.Cells.Item(9, 4).Value = "e2e"
$wsZhCn.Cells.Item(9, 5).Value = "ht"
$wsZhCn.Cells.Item(9, 6).Value = "'False"
$wsZhCn.Cells.Item(9, 7).Value = $xliffZhCn
$wsZhCn.Cells.Item(9, 8).Value = $handoffDateZhCn
$wsZhCn.Cells.Item(9, 9).Value = ""
$wsZhCn.Cells.Item(9, 10).Value = ""
$wsZhCn.Cells.Item(9, 11).Value = ""
$wsZhCn.Cells.Item(9, 12).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(9, 13).Value = ""
$wsZhCn.Cells.Item(9, 14).Value = ""
$wsZhCn.Cells.Item(9, 15).Value = "'True"
$wsZhCn.Cells.Item(9, 16).Value = ""
$wsZhCn.Cells.Item(9, 17).Value = "'False"
$wsZhCn.Cells.Item(9, 18).Value = ""

$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Cells.Item(9, 1).Value = $fileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item(9, 1), $ghBase, "", "", $fileName) | Out-Null
$wsDeDe.Cells.Item(9, 2).Value = ".md"
$wsDeDe.Cells.Item(9, 3).Value = "Ready for handoff"
$wsDeDe.Cells.Item(9, 4).Value = "e2e"
$wsDeDe.Cells.Item(9, 5).Value = "ht"
$wsDeDe.Cells.Item(9, 6).Value = "'False"
$wsDeDe.Cells.Item(9, 7).Value = $xliffDeDe
$wsDeDe.Cells.Item(9, 8).Value = $handoffDateDeDe
$wsDeDe.Cells.Item(9, 9).Value = ""
$wsDeDe.Cells.Item(9, 10).Value = ""
$wsDeDe.Cells.Item(9, 11).Value = ""
$wsDeDe.Cells.Item(9, 12).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(9, 13).Value = ""
$wsDeDe.Cells.Item(9, 14).Value = ""
$wsDeDe.Cells.Item(9, 15).Value = "'True"
$wsDeDe.Cells.Item(9, 16).Value = ""
$wsDeDe.Cells.Item(9, 17).Value = "'False"
$wsDeDe.Cells.Item(9, 18).Value = ""
